$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "INQUINANTE"
$ws.Range("B1").Value = "ANNO"
$ws.Range("C1").Value = "PERCORSO EMS"
$ws.Range("D1").Value = "PERCORSO POLL"
$ws.Range("E1").Value = "PERCORSO GRIGLIA"
$ws.Range("F1").Value = "PERCORSO POPOLAZIONE"
$ws.Range("G1").Value = "PERCORSO KPI"
$ws.Range("H1").Value = "SCELTA GRIGLIA"
$ws.Range("I1").Value = "PERCENTILE POLL MINIMO"
$ws.Range("J1").Value = "PERCENTILE POLL MASSIMO"
$ws.Range("K1").Value = "PASSO POLL"
$ws.Range("L1").Value = "PERCENTILE EMS MINIMO"
$ws.Range("M1").Value = "PERCENTILE EMS MASSIMO"
$ws.Range("N1").Value = "PASSO EMS"
$ws.Range("O1").Value = "MAX or MEAN"
$ws.Range("P1").Value = "PERC or SOGLIA"
$ws.Range("Q1").Value = "VALORE DI SOGLIA"

# --- Data row (row 2) ---
$ws.Range("A2").Value = "PM2.5"
$ws.Range("B2").Value = 19
$ws.Range("C2").Value = "C:\Users\ASUS\Desktop\dati2019_01\AREU-2019_01-CRS32632-datatime.shp"
$ws.Range("D2").Value = "C:\Users\ASUS\Desktop\dati2019_01\dailyPM25-2019_01.shp"
$ws.Range("E2").Value = "C:\Users\ASUS\Desktop\UNI\3° ANNO\PROGETTO\DATI\GriglieGeografiche\LMB3A.shp"
$ws.Range("F2").Value = "C:\Users\ASUS\Desktop\UNI\3° ANNO\PROGETTO\DATI\GriglieGeografiche\Grids_PopData\LMB3A_POP_2018.csv"
$ws.Range("G2").Value = "C:\Users\ASUS\Desktop\KPI"
$ws.Range("H2").Value = "3A"
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 95
$ws.Range("K2").Value = 2.5
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = 95
$ws.Range("N2").Value = 2.5
$ws.Range("O2").Value = "MEAN"
$ws.Range("P2").Value = "PERC"

# --- Ensure header row formatting matches the rest (bold, centered) ---
$ws.Range("L1").Copy()
$ws.Range("M1:Q1").PasteSpecial(-4122)

# --- Column widths ---
$ws.Range("A1:B1").ColumnWidth = 22.109375
$ws.Range("C1").ColumnWidth = 28.5546875
$ws.Range("D1").ColumnWidth = 28.44140625
$ws.Range("E1").ColumnWidth = 29.88671875
$ws.Range("F1").ColumnWidth = 47.77734375
$ws.Range("G1").ColumnWidth = 29.88671875
$ws.Range("H1").ColumnWidth = 25.5546875
$ws.Range("I1:J1").ColumnWidth = 36.33203125
$ws.Range("K1").ColumnWidth = 20.109375
$ws.Range("L1").ColumnWidth = 34.88671875
$ws.Range("M1").ColumnWidth = 39.5546875
$ws.Range("N1").ColumnWidth = 20.5546875
$ws.Range("O1").ColumnWidth = 22
$ws.Range("P1").ColumnWidth = 27.44140625
$ws.Range("Q1").ColumnWidth = 23.109375

# --- View / selection state ---
$ws.Range("Q1").Select()
$excel.ActiveWindow.ScrollColumn = 10
